$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 2919.5925
$ws.Range("I40").Value = 3524.1428
$ws.Range("K40").Value = 3524.1428
$ws.Range("M40").Value = -3349.1428

$ws.Range("H62").Value = 3758
$ws.Range("I62").Value = 2900
$ws.Range("J62").Value = 4902
$ws.Range("K62").Value = 2900
$ws.Range("L62").Value = 4902
$ws.Range("M62").Value = -2276
$ws.Range("N62").Value = -6150

$ws.Range("H64").Value = 2837.913
$ws.Range("I64").Value = 2898.25
$ws.Range("J64").Value = 2805.7334
$ws.Range("K64").Value = 2898.25
$ws.Range("L64").Value = 2805.7334
$ws.Range("M64").Value = -2650.25
$ws.Range("N64").Value = -3301.7334

$ws.Range("H65").Value = 3758
$ws.Range("I65").Value = 2900
$ws.Range("J65").Value = 4902
$ws.Range("K65").Value = 14500
$ws.Range("L65").Value = 24510
$ws.Range("M65").Value = -11380
$ws.Range("N65").Value = -30750

$ws.Range("H67").Value = 2837.913
$ws.Range("I67").Value = 2898.25
$ws.Range("J67").Value = 2805.7334
$ws.Range("K67").Value = 2898.25
$ws.Range("L67").Value = 2805.7334
$ws.Range("M67").Value = -2040.25
$ws.Range("N67").Value = -4521.7334

$ws.Range("H98").Value = 976.5
$ws.Range("I98").Value = 1035
$ws.Range("J98").Value = 884.5714
$ws.Range("K98").Value = 1035
$ws.Range("L98").Value = 884.5714
$ws.Range("M98").Value = 463
$ws.Range("N98").Value = -3880.5714

$ws.Range("H118").Value = 4446.303
$ws.Range("I118").Value = 470
$ws.Range("J118").Value = 8188.706
$ws.Range("K118").Value = 1410
$ws.Range("L118").Value = 24566.118
$ws.Range("M118").Value = 247
$ws.Range("N118").Value = -27880.118

$ws.Range("H122").Value = 976.5
$ws.Range("I122").Value = 1035
$ws.Range("J122").Value = 884.5714
$ws.Range("K122").Value = 3105
$ws.Range("L122").Value = 2653.7142
$ws.Range("M122").Value = -655
$ws.Range("N122").Value = -7553.7142

$ws.Range("H138").Value = 2713.8462
$ws.Range("I138").Value = 1317.9412
$ws.Range("J138").Value = 3208.2292
$ws.Range("K138").Value = 3953.8236
$ws.Range("L138").Value = 9624.687600000001
$ws.Range("M138").Value = 1186.1764
$ws.Range("N138").Value = -19904.6876

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 703
$ws.Range("I2").Value = 656.8333
$ws.Range("J2").Value = 910.75
$ws.Range("K2").Value = 656.8333
$ws.Range("L2").Value = 910.75
$ws.Range("M2").Value = -543.8333
$ws.Range("N2").Value = -1136.75

$ws.Range("H45").Value = 1718.4286
$ws.Range("I45").Value = 1719.0385
$ws.Range("K45").Value = 1719.0385
$ws.Range("M45").Value = -1342.0385

$ws.Range("H88").Value = 142858460
$ws.Range("I88").Value = 1432.6666
$ws.Range("J88").Value = 250001250
$ws.Range("K88").Value = 1432.6666
$ws.Range("L88").Value = 250001250
$ws.Range("M88").Value = -1026.6666
$ws.Range("N88").Value = -250002062

$ws.Range("H91").Value = 142858460
$ws.Range("I91").Value = 1432.6666
$ws.Range("J91").Value = 250001250
$ws.Range("K91").Value = 1432.6666
$ws.Range("L91").Value = 250001250
$ws.Range("M91").Value = -28.66660000000002
$ws.Range("N91").Value = -250004058

$ws.Range("H110").Value = 1290.3462
$ws.Range("I110").Value = 1111.8889
$ws.Range("J110").Value = 1691.875
$ws.Range("K110").Value = 1111.8889
$ws.Range("L110").Value = 1691.875
$ws.Range("M110").Value = 933.1111000000001
$ws.Range("N110").Value = -5781.875

$ws.Range("H116").Value = 703
$ws.Range("I116").Value = 656.8333
$ws.Range("J116").Value = 910.75
$ws.Range("K116").Value = 656.8333
$ws.Range("L116").Value = 910.75
$ws.Range("M116").Value = 1637.1667
$ws.Range("N116").Value = -5498.75

$ws.Range("H122").Value = 1943.0358
$ws.Range("I122").Value = 1504.4783
$ws.Range("K122").Value = 4513.4349
$ws.Range("M122").Value = -2063.4349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 703
$ws.Range("I3").Value = 656.8333
$ws.Range("J3").Value = 910.75
$ws.Range("K3").Value = 656.8333
$ws.Range("L3").Value = 910.75
$ws.Range("M3").Value = -542.8333
$ws.Range("N3").Value = -1138.75

$ws.Range("H86").Value = 15386820
$ws.Range("I86").Value = 25002438
$ws.Range("J86").Value = 1832.8
$ws.Range("K86").Value = 25002438
$ws.Range("L86").Value = 1832.8
$ws.Range("M86").Value = -25001315
$ws.Range("N86").Value = -4078.8

$ws.Range("H89").Value = 15386820
$ws.Range("I89").Value = 25002438
$ws.Range("J89").Value = 1832.8
$ws.Range("K89").Value = 125012190
$ws.Range("L89").Value = 9164
$ws.Range("M89").Value = -125006574
$ws.Range("N89").Value = -20396

$ws.Range("H99").Value = 1616
$ws.Range("I99").Value = 1454.6154
$ws.Range("J99").Value = 1915.7142
$ws.Range("K99").Value = 1454.6154
$ws.Range("L99").Value = 1915.7142
$ws.Range("M99").Value = 43.38460000000009
$ws.Range("N99").Value = -4911.7142

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1691.6666
$ws.Range("I16").Value = 1895
$ws.Range("J16").Value = 1285
$ws.Range("K16").Value = 1895
$ws.Range("L16").Value = 1285
$ws.Range("M16").Value = -1608
$ws.Range("N16").Value = -1859

$ws.Range("H99").Value = 4128.591
$ws.Range("I99").Value = 3991.45
$ws.Range("J99").Value = 5500
$ws.Range("K99").Value = 3991.45
$ws.Range("L99").Value = 5500
$ws.Range("M99").Value = -2493.45
$ws.Range("N99").Value = -8496

$ws.Range("H113").Value = 1691.6666
$ws.Range("I113").Value = 1895
$ws.Range("J113").Value = 1285
$ws.Range("K113").Value = 1895
$ws.Range("L113").Value = 1285
$ws.Range("M113").Value = 275
$ws.Range("N113").Value = -5625

$ws.Range("H122").Value = 973.55554
$ws.Range("I122").Value = 737.5
$ws.Range("J122").Value = 1162.4
$ws.Range("K122").Value = 2212.5
$ws.Range("L122").Value = 3487.2
$ws.Range("M122").Value = 237.5
$ws.Range("N122").Value = -8387.200000000001

$ws.Range("H126").Value = 4128.591
$ws.Range("I126").Value = 3991.45
$ws.Range("J126").Value = 5500
$ws.Range("K126").Value = 11974.35
$ws.Range("L126").Value = 16500
$ws.Range("M126").Value = -9504.349999999999
$ws.Range("N126").Value = -21440

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 992.8570999999999
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 992.8570999999999
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").Value = 2978.5713
$ws.Range("N80").Value = -4850.5713

$ws.Range("H83").Value = 992.8570999999999
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 992.8570999999999
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").Value = 8935.713899999999
$ws.Range("N83").Value = -18295.7139

$ws.Range("H113").Value = 863.871
$ws.Range("I113").Value = 422.84616
$ws.Range("J113").Value = 1182.3889
$ws.Range("K113").Value = 1268.53848
$ws.Range("L113").Value = 3547.1667
$ws.Range("M113").Value = 901.4615200000001
$ws.Range("N113").Value = -7887.1667

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H74").Value = 30000
$ws.Range("J74").Value = 30000
$ws.Range("L74").Value = 30000
$ws.Range("N74").Value = -31872

$ws.Range("H75").Value = 24000
$ws.Range("J75").Value = 24000
$ws.Range("L75").Value = 24000
$ws.Range("N75").Value = -25748

$ws.Range("H77").Value = 30000
$ws.Range("J77").Value = 30000
$ws.Range("L77").Value = 90000
$ws.Range("N77").Value = -99360

$ws.Range("H78").Value = 24000
$ws.Range("J78").Value = 24000
$ws.Range("L78").Value = 72000
$ws.Range("N78").Value = -80736

$ws.Range("H86").Value = 30095.334
$ws.Range("J86").Value = 30095.334
$ws.Range("L86").Value = 30095.334
$ws.Range("N86").Value = -32467.334

$ws.Range("H89").Value = 30095.334
$ws.Range("J89").Value = 30095.334
$ws.Range("L89").Value = 90286.00199999999
$ws.Range("N89").Value = -102142.002

$ws.Range("H102").Value = 918.6
$ws.Range("I102").Value = 888.75
$ws.Range("K102").Value = 888.75
$ws.Range("M102").Value = 733.25

$ws.Range("H113").Value = 1347.5625
$ws.Range("I113").Value = 1023.7273
$ws.Range("J113").Value = 2060
$ws.Range("K113").Value = 1023.7273
$ws.Range("L113").Value = 2060
$ws.Range("M113").Value = 1146.2727
$ws.Range("N113").Value = -6400

$ws.Range("H122").Value = 5678.5
$ws.Range("I122").Value = 5678.5
$ws.Range("K122").Value = 17035.5
$ws.Range("M122").Value = -14585.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1548.7142
$ws.Range("I40").Value = 1508.7858
$ws.Range("J40").Value = 1628.5714
$ws.Range("K40").Value = 1508.7858
$ws.Range("L40").Value = 1628.5714
$ws.Range("M40").Value = -1372.7858
$ws.Range("N40").Value = -1900.5714

$ws.Range("H61").Value = 1796
$ws.Range("I61").Value = 1850.8125
$ws.Range("J61").Value = 1670.7142
$ws.Range("K61").Value = 1850.8125
$ws.Range("L61").Value = 1670.7142
$ws.Range("M61").Value = -1648.8125
$ws.Range("N61").Value = -2074.7142

$ws.Range("H113").Value = 1796
$ws.Range("I113").Value = 1850.8125
$ws.Range("J113").Value = 1670.7142
$ws.Range("K113").Value = 1850.8125
$ws.Range("L113").Value = 1670.7142
$ws.Range("M113").Value = 319.1875
$ws.Range("N113").Value = -6010.7142

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 334.57144
$ws.Range("I113").Value = 266.45456
$ws.Range("J113").Value = 584.3333
$ws.Range("K113").Value = 799.36368
$ws.Range("L113").Value = 1752.9999
$ws.Range("M113").Value = 1370.63632
$ws.Range("N113").Value = -6092.9999

$ws.Range("H122").Value = 1101.2059
$ws.Range("I122").Value = 868
$ws.Range("J122").Value = 1396.6
$ws.Range("K122").Value = 2604
$ws.Range("L122").Value = 4189.799999999999
$ws.Range("M122").Value = -154
$ws.Range("N122").Value = -9089.799999999999

$ws.Range("H126").Value = 2977.3157
$ws.Range("I126").Value = 3511.7144
$ws.Range("J126").Value = 1481
$ws.Range("K126").Value = 10535.1432
$ws.Range("L126").Value = 4443
$ws.Range("M126").Value = -8065.143199999999
$ws.Range("N126").Value = -9383
